$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E5").Value = 1
$ws.Range("E15").Value = 101
$ws.Range("E18").Value = 53
$ws.Range("E42").Value = 21
$ws.Range("E44").Value = 17
$ws.Range("E49").Value = 35
$ws.Range("E57").Value = 7
$ws.Range("F67").Value = 10
$ws.Range("H67").Value = 10
$ws.Range("E77").Value = 27
$ws.Range("E82").Value = 4
